# Rename the (only) worksheet from "Sheet1" to "Sheet 1"
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Name = "Sheet 1"

# Move the active selection from E11 to C14
$ws.Range("C14").Select()
